# Updates Kujata_Profits market-profit values across multiple leve sheets
# (scheduled data refresh of currentAveragePrice / LeveProfit columns).
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Cells.Item(62, 8).Value = 15875272
$ws.Cells.Item(62, 9).Value = 27778978
$ws.Cells.Item(62, 11).Value = 27778978
$ws.Cells.Item(62, 13).Value = -27778354
# Row 65
$ws.Cells.Item(65, 8).Value = 15875272
$ws.Cells.Item(65, 9).Value = 27778978
$ws.Cells.Item(65, 11).Value = 138894890
$ws.Cells.Item(65, 13).Value = -138891770
# Row 86
$ws.Cells.Item(86, 8).Value = 4099.5
$ws.Cells.Item(86, 9).Value = 4199.3335
$ws.Cells.Item(86, 10).Value = 3800
$ws.Cells.Item(86, 11).Value = 4199.3335
$ws.Cells.Item(86, 12).Value = 3800
$ws.Cells.Item(86, 13).Value = -3076.3335
$ws.Cells.Item(86, 14).Value = -6046
# Row 87
$ws.Cells.Item(87, 8).Value = 28848.75
$ws.Cells.Item(87, 10).Value = 28848.75
$ws.Cells.Item(87, 12).Value = 28848.75
$ws.Cells.Item(87, 14).Value = -31344.75
# Row 88
$ws.Cells.Item(88, 8).Value = 823612.8
$ws.Cells.Item(88, 9).Value = 485.14285
$ws.Cells.Item(88, 10).Value = 1543849.5
$ws.Cells.Item(88, 11).Value = 485.14285
$ws.Cells.Item(88, 12).Value = 1543849.5
$ws.Cells.Item(88, 13).Value = -79.14285000000001
$ws.Cells.Item(88, 14).Value = -1544661.5
# Row 89
$ws.Cells.Item(89, 8).Value = 4099.5
$ws.Cells.Item(89, 9).Value = 4199.3335
$ws.Cells.Item(89, 10).Value = 3800
$ws.Cells.Item(89, 11).Value = 20996.6675
$ws.Cells.Item(89, 12).Value = 19000
$ws.Cells.Item(89, 13).Value = -15380.6675
$ws.Cells.Item(89, 14).Value = -30232
# Row 90
$ws.Cells.Item(90, 8).Value = 28848.75
$ws.Cells.Item(90, 10).Value = 28848.75
$ws.Cells.Item(90, 12).Value = 86546.25
$ws.Cells.Item(90, 14).Value = -99026.25
# Row 91
$ws.Cells.Item(91, 8).Value = 823612.8
$ws.Cells.Item(91, 9).Value = 485.14285
$ws.Cells.Item(91, 10).Value = 1543849.5
$ws.Cells.Item(91, 11).Value = 485.14285
$ws.Cells.Item(91, 12).Value = 1543849.5
$ws.Cells.Item(91, 13).Value = 918.85715
$ws.Cells.Item(91, 14).Value = -1546657.5
# Row 106
$ws.Cells.Item(106, 8).Value = 7249.4
$ws.Cells.Item(106, 9).Value = 7473.0527
$ws.Cells.Item(106, 11).Value = 7473.0527
$ws.Cells.Item(106, 13).Value = -6842.0527
# Row 133
$ws.Cells.Item(133, 8).Value = 34089
$ws.Cells.Item(133, 10).Value = 34089
$ws.Cells.Item(133, 12).Value = 34089
$ws.Cells.Item(133, 14).Value = -44209

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Cells.Item(45, 8).Value = 1206.75
$ws.Cells.Item(45, 9).Value = 1060.125
$ws.Cells.Item(45, 11).Value = 1060.125
$ws.Cells.Item(45, 13).Value = -683.125
# Row 97
$ws.Cells.Item(97, 8).Value = 466.58334
$ws.Cells.Item(97, 9).Value = 466.58334
$ws.Cells.Item(97, 11).Value = 466.58334
$ws.Cells.Item(97, 13).Value = 29.41665999999998
# Row 122
$ws.Cells.Item(122, 8).Value = 1909.0869
$ws.Cells.Item(122, 9).Value = 1743.2
$ws.Cells.Item(122, 10).Value = 2220.125
$ws.Cells.Item(122, 11).Value = 5229.6
$ws.Cells.Item(122, 12).Value = 6660.375
$ws.Cells.Item(122, 13).Value = -2779.6
$ws.Cells.Item(122, 14).Value = -11560.375
# Row 132
$ws.Cells.Item(132, 8).Value = 2949.9707
$ws.Cells.Item(132, 9).Value = 2591.611
$ws.Cells.Item(132, 10).Value = 3353.125
$ws.Cells.Item(132, 11).Value = 7774.833
$ws.Cells.Item(132, 12).Value = 10059.375
$ws.Cells.Item(132, 13).Value = -5244.833
$ws.Cells.Item(132, 14).Value = -15119.375

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Cells.Item(94, 8).Value = 9615726
$ws.Cells.Item(94, 9).Value = 10869915
$ws.Cells.Item(94, 10).Value = 279.33334
$ws.Cells.Item(94, 11).Value = 10869915
$ws.Cells.Item(94, 12).Value = 279.33334
$ws.Cells.Item(94, 13).Value = -10869464
$ws.Cells.Item(94, 14).Value = -1181.33334
# Row 99
$ws.Cells.Item(99, 8).Value = 71429570
$ws.Cells.Item(99, 9).Value = 90910050
$ws.Cells.Item(99, 10).Value = 1166.6666
$ws.Cells.Item(99, 11).Value = 90910050
$ws.Cells.Item(99, 12).Value = 1166.6666
$ws.Cells.Item(99, 13).Value = -90908552
$ws.Cells.Item(99, 14).Value = -4162.6666
# Row 134
$ws.Cells.Item(134, 8).Value = 1277.3572
$ws.Cells.Item(134, 9).Value = 1073.5834
$ws.Cells.Item(134, 11).Value = 3220.7502
$ws.Cells.Item(134, 13).Value = -685.7501999999999

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
# Row 62
$ws.Cells.Item(62, 9).Value = 8234.286
$ws.Cells.Item(62, 10).Value = 50001724
$ws.Cells.Item(62, 11).Value = 8234.286
$ws.Cells.Item(62, 12).Value = 50001724
$ws.Cells.Item(62, 13).Value = -7610.286
$ws.Cells.Item(62, 14).Value = -50002972
# Row 65
$ws.Cells.Item(65, 9).Value = 8234.286
$ws.Cells.Item(65, 10).Value = 50001724
$ws.Cells.Item(65, 11).Value = 41171.43
$ws.Cells.Item(65, 12).Value = 250008620
$ws.Cells.Item(65, 13).Value = -38051.43
$ws.Cells.Item(65, 14).Value = -250014860
# Row 99
$ws.Cells.Item(99, 8).Value = 1652.2858
$ws.Cells.Item(99, 9).Value = 1581.3334
$ws.Cells.Item(99, 11).Value = 1581.3334
$ws.Cells.Item(99, 13).Value = -83.33339999999998
# Row 107
$ws.Cells.Item(107, 8).Value = 560.7619
$ws.Cells.Item(107, 9).Value = 428
$ws.Cells.Item(107, 10).Value = 1125
$ws.Cells.Item(107, 11).Value = 428
$ws.Cells.Item(107, 12).Value = 1125
$ws.Cells.Item(107, 13).Value = 1492
$ws.Cells.Item(107, 14).Value = -4965
# Row 122
$ws.Cells.Item(122, 8).Value = 737.4706
$ws.Cells.Item(122, 9).Value = 751.6429000000001
$ws.Cells.Item(122, 11).Value = 2254.9287
$ws.Cells.Item(122, 13).Value = 195.0712999999996
# Row 126
$ws.Cells.Item(126, 8).Value = 1652.2858
$ws.Cells.Item(126, 9).Value = 1581.3334
$ws.Cells.Item(126, 11).Value = 4744.0002
$ws.Cells.Item(126, 13).Value = -2274.0002
# Row 134
$ws.Cells.Item(134, 8).Value = 20835318
$ws.Cells.Item(134, 9).Value = 1944
$ws.Cells.Item(134, 10).Value = 50002040
$ws.Cells.Item(134, 11).Value = 5832
$ws.Cells.Item(134, 12).Value = 150006120
$ws.Cells.Item(134, 13).Value = -3297
$ws.Cells.Item(134, 14).Value = -150011190

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Cells.Item(131, 8).Value = 32262124
$ws.Cells.Item(131, 10).Value = 6028.25
$ws.Cells.Item(131, 12).Value = 18084.75
$ws.Cells.Item(131, 14).Value = -28164.75

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
# Row 36
$ws.Cells.Item(36, 8).Value = 1458.5
$ws.Cells.Item(36, 9).Value = 1458.5
$ws.Cells.Item(36, 11).Value = 1458.5
$ws.Cells.Item(36, 13).Value = -973.5
# Row 46
$ws.Cells.Item(46, 8).Value = 13000
$ws.Cells.Item(46, 10).Value = 18000
$ws.Cells.Item(46, 12).Value = 18000
$ws.Cells.Item(46, 14).Value = -18312
# Row 70
$ws.Cells.Item(70, 8).Value = 34618636
$ws.Cells.Item(70, 9).Value = 31253356
$ws.Cells.Item(70, 11).Value = 31253356
$ws.Cells.Item(70, 13).Value = -31253086
# Row 73
$ws.Cells.Item(73, 8).Value = 34618636
$ws.Cells.Item(73, 9).Value = 31253356
$ws.Cells.Item(73, 11).Value = 31253356
$ws.Cells.Item(73, 13).Value = -31252420
# Row 102
$ws.Cells.Item(102, 8).Value = 1304.881
$ws.Cells.Item(102, 9).Value = 1211.8485
$ws.Cells.Item(102, 11).Value = 1211.8485
$ws.Cells.Item(102, 13).Value = 410.1514999999999
# Row 126
$ws.Cells.Item(126, 8).Value = 2069.75
$ws.Cells.Item(126, 9).Value = 1729.5454
$ws.Cells.Item(126, 11).Value = 5188.6362
$ws.Cells.Item(126, 13).Value = -2718.6362
# Row 134
$ws.Cells.Item(134, 8).Value = 29710.166
$ws.Cells.Item(134, 10).Value = 29710.166
$ws.Cells.Item(134, 12).Value = 89130.49800000001
$ws.Cells.Item(134, 14).Value = -94200.49800000001

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Cells.Item(68, 8).Value = 1283
$ws.Cells.Item(68, 10).Value = 0
$ws.Cells.Item(68, 12).Value = 0
$ws.Cells.Item(68, 14).Value = ""
# Row 71
$ws.Cells.Item(71, 8).Value = 1283
$ws.Cells.Item(71, 10).Value = 0
$ws.Cells.Item(71, 12).Value = 0
$ws.Cells.Item(71, 14).Value = ""
# Row 133
$ws.Cells.Item(133, 8).Value = 45936
$ws.Cells.Item(133, 10).Value = 45936
$ws.Cells.Item(133, 12).Value = 45936
$ws.Cells.Item(133, 14).Value = -50996

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
# Row 26
$ws.Cells.Item(26, 8).Value = 3000
$ws.Cells.Item(26, 9).Value = 3000
$ws.Cells.Item(26, 11).Value = 3000
$ws.Cells.Item(26, 13).Value = -2707
# Row 28
$ws.Cells.Item(28, 8).Value = 0
$ws.Cells.Item(28, 10).Value = 0
$ws.Cells.Item(28, 12).Value = 0
$ws.Cells.Item(28, 14).Value = ""
# Row 32
$ws.Cells.Item(32, 8).Value = 2800
$ws.Cells.Item(32, 9).Value = 2800
$ws.Cells.Item(32, 11).Value = 2800
$ws.Cells.Item(32, 13).Value = -2483
# Row 33
$ws.Cells.Item(33, 8).Value = 15000
$ws.Cells.Item(33, 10).Value = 15000
$ws.Cells.Item(33, 12).Value = 15000
$ws.Cells.Item(33, 14).Value = -15500
# Row 36
$ws.Cells.Item(36, 8).Value = 15000
$ws.Cells.Item(36, 10).Value = 15000
$ws.Cells.Item(36, 12).Value = 15000
$ws.Cells.Item(36, 14).Value = -15500
# Row 43
$ws.Cells.Item(43, 8).Value = 1250
$ws.Cells.Item(43, 9).Value = 500
$ws.Cells.Item(43, 10).Value = 2000
$ws.Cells.Item(43, 11).Value = 500
$ws.Cells.Item(43, 12).Value = 2000
$ws.Cells.Item(43, 13).Value = -351
$ws.Cells.Item(43, 14).Value = -2298
# Row 62
$ws.Cells.Item(62, 8).Value = 125003250
$ws.Cells.Item(62, 9).Value = 166669330
$ws.Cells.Item(62, 10).Value = 5000
$ws.Cells.Item(62, 11).Value = 166669330
$ws.Cells.Item(62, 12).Value = 5000
$ws.Cells.Item(62, 13).Value = -166668706
$ws.Cells.Item(62, 14).Value = -6248
# Row 65
$ws.Cells.Item(65, 8).Value = 125003250
$ws.Cells.Item(65, 9).Value = 166669330
$ws.Cells.Item(65, 10).Value = 5000
$ws.Cells.Item(65, 11).Value = 833346650
$ws.Cells.Item(65, 12).Value = 25000
$ws.Cells.Item(65, 13).Value = -833343530
$ws.Cells.Item(65, 14).Value = -31240
# Row 122
$ws.Cells.Item(122, 8).Value = 12501914
$ws.Cells.Item(122, 9).Value = 13890838
$ws.Cells.Item(122, 11).Value = 41672514
$ws.Cells.Item(122, 13).Value = -41670064
# Row 126
$ws.Cells.Item(126, 8).Value = 50001170
$ws.Cells.Item(126, 9).Value = 83333864
$ws.Cells.Item(126, 10).Value = 2126.75
$ws.Cells.Item(126, 11).Value = 250001592
$ws.Cells.Item(126, 12).Value = 6380.25
$ws.Cells.Item(126, 13).Value = -249999122
$ws.Cells.Item(126, 14).Value = -11320.25
